# "MANIJA FIJA" price list - refresh date stamp and unit prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date stamp in A1 (serial date number).
$ws.Range("A1").Value = 45436

# Updated prices in column D for rows 34-37.
$ws.Range("D34").Value = 206.846
$ws.Range("D35").Value = 293.075
$ws.Range("D36").Value = 396.525
$ws.Range("D37").Value = 653.796
